# 6-Jul-2021, end of day update.
# Applies the day's petty-cash entries to the "Sheet1" daily ledger.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 3: top-up the existing cash advance formula with an extra 260,000 ---
$ws.Range("D3").Formula = "=60000+260000"

# --- Row 4: add the 76,848,000 receipt to the existing formula total ---
$ws.Range("C4").Formula = "=1000000+2681000+1545000+76848000"

# --- Row 7: PLN (electricity) payment ---
$ws.Range("B7").Value = "PLN - Astar 214"
$ws.Range("D7").Value = 102500

# --- Row 8: cash/retail sales receipt ---
$ws.Range("B8").Value = "SALES - cash/retail"
$ws.Range("C8").Formula = "=5944275+82502725-76848000"

# --- Row 9: cash overage ---
$ws.Range("B9").Value = "SELISIH - lebih"
$ws.Range("C9").Value = 20000

# --- Row 10: bank deposit ---
$ws.Range("B10").Value = "SETOR KE BANK"
$ws.Range("D10").Value = 82000000

# --- Row 11: new day, 7-Jul-2021, meal allowance (Wages Expense) ---
$ws.Range("A11").Value = 44383
$ws.Range("B11").Value = "Wages Expense"
$ws.Range("D11").Formula = "=60000"

# --- Row 12: tax payment - P.Tata ---
$ws.Range("B12").Value = "TAX - P.Tata"
$ws.Range("D12").Value = 200000

# --- Row 13: tax payment - Iuran ARIESTA ---
$ws.Range("B13").Value = "TAX - Iuran ARIESTA"
$ws.Range("D13").Value = 660000

# --- Row 14: BCA transfer ---
$ws.Range("B14").Value = "TRANSFER BCA"
$ws.Range("D14").Formula = "=4500000+432000+632000"

# --- Row 15: freight out ---
$ws.Range("B15").Value = "FREIGHT - OUT"
$ws.Range("D15").Value = 8000

# --- Recalculate so the running E-column balances update ---
$excel.Calculate()

# --- Update the frozen-pane scroll position / active selection to reflect
#     where the user was working at end-of-day (row 10 at top, D31 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("D31").Select() | Out-Null

$wb.Save()
